$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B28 was stored as text "081245789"; it becomes a real number (81245789).
$ws.Range("B28").Value = 81245789

# New rows of order data appended below the existing data (rows 29-33).
# Columns: Nama, Nomor Telepon, Alamat, Jumlah, Metode, Total,
#          Tanggal Pengiriman, Waktu Pengiriman, Catatan, Status, Bukti Pembayaran
$rows = @(
    @("Alber",  "1243569",     "Pandan Alas",     15, "Transfer Bank", 375000,  "2025-05-18", "20:07"),
    @("aurora", "47859",       "Jalan Pete Raya",  70, "QRIS",          1750000, "2025-05-18", "20:07"),
    @("Natavia","1245",        "GG Cempaka",       70, "Transfer Bank", 1750000, "2025-05-18", "20:15"),
    @("Larisa", "081315130649","Jl Padanaran",      1, "Transfer Bank", 25000,   "2025-05-18", "20:29"),
    @("Aurora", "1478956231",  "GG. Pete",          1, "QRIS",          25000,   "2025-05-18", "20:32")
)

$r = 29
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Phone numbers must stay textual even though they look numeric - force
    # text storage, then drop the temporary "@" format override so the cell
    # is left without any explicit style (matching plain data rows).
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).ClearFormats()

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # The "Tanggal Pengiriman" column holds a plain yyyy-mm-dd text label,
    # not a real date - same text-forcing trick as column B.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).ClearFormats()

    $ws.Cells.Item($r, 8).Value = $row[7]

    # Catatan / Bukti Pembayaran are blank for every new order.
    $ws.Cells.Item($r, 9).Value = ""

    $ws.Cells.Item($r, 10).Value = "Diproses"

    $ws.Cells.Item($r, 11).Value = ""

    $r = $r + 1
}
